$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings
# like "26.341.69" / "1.009" are stored verbatim as text, matching
# the inline-string cells in the source workbook (no float coercion).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.341.69"
$ws.Range("E2").Value = "  +0.26%  "

$ws.Range("D3").Value = "1.686.76"
$ws.Range("E3").Value = "  +1.07%  "

$ws.Range("E4").Value = "  +0.36%  "

$ws.Range("D5").Value = "218.54"
$ws.Range("E5").Value = "  +0.43%  "

$ws.Range("D6").Value = "0.5257"
$ws.Range("E6").Value = "  +3.51%  "

$ws.Range("D7").Value = "1.009"
$ws.Range("E7").Value = "  +0.41%  "

$ws.Range("D8").Value = "0.2704"
$ws.Range("E8").Value = "  +1.90%  "

$ws.Range("B9").Value = "Solana"
$ws.Range("C9").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D9").Value = "22.06"
$ws.Range("E9").Value = "  +2.68%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.06417"
$ws.Range("E10").Value = "  +1.01%  "

$ws.Range("D11").Value = "0.07497"
$ws.Range("E11").Value = "  +1.89%  "

$ws.Range("D12").Value = "1.726.93"
$ws.Range("E12").Value = "  +3.39%  "

$ws.Range("D13").Value = "4.573"
$ws.Range("E13").Value = "  +0.62%  "

$ws.Range("D14").Value = "0.5824"
$ws.Range("E14").Value = "  +0.29%  "

$ws.Range("D15").Value = "0.000008499"
$ws.Range("E15").Value = "  -0.16%  "

$ws.Range("D16").Value = "64.51"
$ws.Range("E16").Value = "  -0.32%  "

$ws.Range("D17").Value = "26.387.73"
$ws.Range("E17").Value = "  +0.24%  "

$ws.Range("D18").Value = "4.936"
$ws.Range("E18").Value = "  +0.26%  "

$ws.Range("D19").Value = "1.008"
$ws.Range("E19").Value = "  +0.21%  "

$ws.Range("D20").Value = "10.89"
$ws.Range("E20").Value = "  +0.52%  "

$ws.Range("D21").Value = "189.58"
$ws.Range("E21").Value = "  +0.48%  "

$ws.Range("D22").Value = "6.212"
$ws.Range("E22").Value = "  +0.38%  "

$ws.Range("E23").Value = "  +0.23%  "

$ws.Range("D24").Value = "144.61"
$ws.Range("E24").Value = "  +0.95%  "

$ws.Range("D25").Value = "7.724"
$ws.Range("E25").Value = "  +0.70%  "

$ws.Range("D26").Value = "0.1232"
$ws.Range("E26").Value = "  +4.96%  "

$ws.Range("D27").Value = "15.84"
$ws.Range("E27").Value = "  +1.35%  "

$ws.Range("D28").Value = "0.06656"
$ws.Range("E28").Value = "  +13.80%  "

$ws.Range("D29").Value = "1.351"
$ws.Range("E29").Value = "  +6.42%  "

$ws.Range("E30").Value = "  +0.49%  "

$ws.Range("D31").Value = "3.578"
$ws.Range("E31").Value = "  +1.93%  "

$ws.Range("D32").Value = "3.572"
$ws.Range("E32").Value = "  +1.27%  "

$ws.Range("D33").Value = "1.661"
$ws.Range("E33").Value = "  +1.39%  "

$ws.Range("D34").Value = "1.027"
$ws.Range("E34").Value = "  +1.56%  "

$ws.Range("D35").Value = "0.6235"
$ws.Range("E35").Value = "  +4.04%  "

$ws.Range("D36").Value = "2.399"
$ws.Range("E36").Value = "  +1.91%  "

$ws.Range("D37").Value = "2.703"
$ws.Range("E37").Value = "  +2.18%  "

$ws.Range("D38").Value = "6.382"
$ws.Range("E38").Value = "  +6.27%  "

$ws.Range("D39").Value = "1.113.13"
$ws.Range("E39").Value = "  +3.75%  "

$ws.Range("D40").Value = "0.01620"
$ws.Range("E40").Value = "  +0.52%  "

$ws.Range("D41").Value = "0.8841"
$ws.Range("E41").Value = "  +2.27%  "

$ws.Range("E42").Value = "  +0.96%  "

$ws.Range("D43").Value = "101.09"
$ws.Range("E43").Value = "  +1.51%  "

$ws.Range("D44").Value = "1.837.22"
$ws.Range("E44").Value = "  +1.19%  "

$ws.Range("D45").Value = "0.00000000111"
$ws.Range("E45").Value = "  +0.77%  "

$ws.Range("D46").Value = "56.98"
$ws.Range("E46").Value = "  +2.22%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "8.184"
$ws.Range("E47").Value = "  +1.25%  "

$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "1.009"
$ws.Range("E48").Value = "  +0.52%  "

$ws.Range("D49").Value = "0.05269"
$ws.Range("E49").Value = "  +1.76%  "

$ws.Range("E50").Value = "  +0.33%  "

$ws.Range("D51").Value = "6.077"
$ws.Range("E51").Value = "  +3.84%  "
